$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.040.26'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.631.46'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.07'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.84'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.649'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +3.59%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.121'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -4.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.81'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.387'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.63'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.51%  '
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.113.57'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000185'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -6.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.914.38'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.659.51'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.19'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.67'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.57'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '346.67'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.79%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.77'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.75'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000109'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.34'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '580.58'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +9.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.58'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.97'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.81%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.161'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.96%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.08'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.71'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.50'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.28'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.411'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.00'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.68%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.89'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '152.24'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '41.92'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '158.30'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.36'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.99'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.20'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0597'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.33%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.103'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.633'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0253'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.53%  '
